# Update crypto price/volume table cells per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D values that are plain decimals (e.g. "580.45") would be
# auto-coerced to numbers (and lose exact text, e.g. trailing zeros)
# by a plain .Value assignment, so force those specific cells to text
# first to preserve the original string formatting used throughout
# this sheet (prices are stored as text, not numbers).
$textCells = @(
  "D5",
  "D6",
  "D8",
  "D11",
  "D12",
  "D14",
  "D16",
  "D20",
  "D21",
  "D22",
  "D23",
  "D24",
  "D25",
  "D27",
  "D28",
  "D29",
  "D30",
  "D31",
  "D32",
  "D34",
  "D35",
  "D38",
  "D40",
  "D41",
  "D43",
  "D44",
  "D46",
  "D48",
  "D50",
  "D51"
)
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.741.96"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").Value = "3.304.91"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "580.45"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "179.11"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("D9").Value = "3.297.60"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").Value = "0.577"
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").Value = "45.97"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  +3.85%  "
$ws.Range("D14").Value = "694.60"
$ws.Range("E14").Value = "  +13.48%  "
$ws.Range("D15").Value = "3.832.94"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "8.42"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "67.798.50"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "3.297.79"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").Value = "17.52"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "10.86"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "0.898"
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").Value = "17.03"
$ws.Range("E23").Value = "  -6.45%  "
$ws.Range("D24").Value = "5.20"
$ws.Range("E24").Value = "  +5.77%  "
$ws.Range("D25").Value = "98.42"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").Value = "2.76"
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").Value = "9.42"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "33.10"
$ws.Range("E29").Value = "  +8.23%  "
$ws.Range("D30").Value = "8.53"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("D31").Value = "6.82"
$ws.Range("E31").Value = "  +6.21%  "
$ws.Range("D32").Value = "584.89"
$ws.Range("E32").Value = "  +7.46%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "3.891.95"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "10.92"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").Value = "0.105"
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  -7.92%  "
$ws.Range("D38").Value = "55.36"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("E39").Value = "  +2.72%  "
$ws.Range("D40").Value = "3.22"
$ws.Range("E40").Value = "  +2.82%  "
$ws.Range("D41").Value = "2.65"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0692"
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "32.48"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "3.36"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").Value = "0.0414"
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("D48").Value = "1.39"
$ws.Range("E48").Value = "  +10.00%  "
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "2.56"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("D51").Value = "128.62"
$ws.Range("E51").Value = "  +0.64%  "
